# Re-order the test-case blocks in the "Test-Cases" sheet: the
# TestScenario_1.TestCase_1 ("New Account") block, previously the last
# block (rows 11-14), moves to become the first block (rows 2-5); every
# other block shifts down by 4 rows accordingly. Rewritten explicitly as
# the full A2:K14 body so the result matches exactly regardless of the
# previous layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 13,11
$data[0,0] = "TestScenario_1"
$data[0,1] = "TestScenario_1.TestCase_1"
$data[0,2] = "New Account"
$data[0,3] = "User Needs to Login to Salesforce, from the browser with correct credentials"
$data[0,4] = ""
$data[0,5] = "Step 1"
$data[0,6] = "Click on the Account tab, and click on New button"
$data[0,7] = "User should be navigated to the New  Account Page"
$data[0,8] = "Approved"
$data[0,9] = ""
$data[0,10] = ""
$data[1,0] = ""
$data[1,1] = ""
$data[1,2] = ""
$data[1,3] = ""
$data[1,4] = "Valid value for required field Account Name  @ Valid value for required field Account Name, value should be  Test for Workflow Process."
$data[1,5] = "Step 2"
$data[1,6] = "Input valid value in the  Account Name field."
$data[1,7] = "User should be able to input value for the Account Name field."
$data[1,8] = ""
$data[1,9] = ""
$data[1,10] = ""
$data[2,0] = ""
$data[2,1] = ""
$data[2,2] = ""
$data[2,3] = ""
$data[2,4] = "Valid value for required field  "
$data[2,5] = "Step 3"
$data[2,6] = "Input valid value in the   field."
$data[2,7] = "User should be able to input value for the  field."
$data[2,8] = ""
$data[2,9] = ""
$data[2,10] = ""
$data[3,0] = ""
$data[3,1] = ""
$data[3,2] = ""
$data[3,3] = ""
$data[3,4] = ""
$data[3,5] = "Step 4"
$data[3,6] = "Click on Save button to save Account with fields"
$data[3,7] = "User should be able to validate that a New Account is created"
$data[3,8] = ""
$data[3,9] = ""
$data[3,10] = ""
$data[4,0] = "TestScenario_1"
$data[4,1] = "TestScenario_2.TestCase_1"
$data[4,2] = "View Account"
$data[4,3] = "User Needs to Login to Salesforce, from the browser with correct credentials"
$data[4,4] = ""
$data[4,5] = "Step 1"
$data[4,6] = "Click on the Account tab,  and select a Account "
$data[4,7] = "User should be navigated to the Account Page"
$data[4,8] = ""
$data[4,9] = ""
$data[4,10] = ""
$data[5,0] = ""
$data[5,1] = ""
$data[5,2] = ""
$data[5,3] = ""
$data[5,4] = ""
$data[5,5] = "Step 2"
$data[5,6] = "Click on Account name to View the Details"
$data[5,7] = "User should be able to view the Account Details"
$data[5,8] = ""
$data[5,9] = ""
$data[5,10] = ""
$data[6,0] = "TestScenario_2"
$data[6,1] = "TestScenario_3.TestCase_1"
$data[6,2] = "Edit Account"
$data[6,3] = "User Needs to Login to Salesforce, from the browser with correct credentials"
$data[6,4] = ""
$data[6,5] = "Step 1"
$data[6,6] = "Click on the Account tab,  and click on existing Account to modify"
$data[6,7] = "User is navigated to the Account Details page"
$data[6,8] = ""
$data[6,9] = ""
$data[6,10] = ""
$data[7,0] = ""
$data[7,1] = ""
$data[7,2] = ""
$data[7,3] = ""
$data[7,4] = "Valid value for required field Account Name  @ Valid value for required field Account Name, value should be  Test for Workflow Process."
$data[7,5] = "Step 2"
$data[7,6] = "Input valid value in the  Account Name field."
$data[7,7] = "User should be able to input value for the Account Name field."
$data[7,8] = ""
$data[7,9] = ""
$data[7,10] = ""
$data[8,0] = ""
$data[8,1] = ""
$data[8,2] = ""
$data[8,3] = ""
$data[8,4] = "Valid value for required field  "
$data[8,5] = "Step 3"
$data[8,6] = "Input valid value in the   field."
$data[8,7] = "User should be able to input value for the  field."
$data[8,8] = ""
$data[8,9] = ""
$data[8,10] = ""
$data[9,0] = ""
$data[9,1] = ""
$data[9,2] = ""
$data[9,3] = ""
$data[9,4] = ""
$data[9,5] = "Step 4"
$data[9,6] = "Click on Save button to save Account with fields"
$data[9,7] = "User should be able to validate that the Account is edited"
$data[9,8] = ""
$data[9,9] = ""
$data[9,10] = ""
$data[10,0] = "TestScenario_3"
$data[10,1] = "TestScenario_4.TestCase_1"
$data[10,2] = "Delete Account"
$data[10,3] = "User Needs to Login to Salesforce, from the browser with correct credentials"
$data[10,4] = ""
$data[10,5] = "Step 1"
$data[10,6] = "Click on the Account tab,  and select the existing  Account to delete"
$data[10,7] = "User is navigated to the Account Details page"
$data[10,8] = ""
$data[10,9] = ""
$data[10,10] = ""
$data[11,0] = ""
$data[11,1] = ""
$data[11,2] = ""
$data[11,3] = ""
$data[11,4] = ""
$data[11,5] = "Step 2"
$data[11,6] = "Click on to the Delete to Delete the Account"
$data[11,7] = "User should be able to validate that a pop-up is displayed asking for confirmation to delete the Account"
$data[11,8] = ""
$data[11,9] = ""
$data[11,10] = ""
$data[12,0] = ""
$data[12,1] = ""
$data[12,2] = ""
$data[12,3] = ""
$data[12,4] = ""
$data[12,5] = "Step 3"
$data[12,6] = "Click on Confirm / OK to delete the  Account"
$data[12,7] = "User should be able to validate the Account is deleted"
$data[12,8] = ""
$data[12,9] = ""
$data[12,10] = ""
$ws.Range("A2:K14").Value = $data
